$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Paste the Neo4j/Cypher MATCH query used to build this test-case workbook
# into A2 (single-quoted PS string so the back-ticked column aliases and
# embedded single quotes survive literally).
$query = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Lung cancer, NOS''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("A2").Value = $query

# Row 2 grows taller to fit the wrapped query text
$ws.Rows.Item(2).RowHeight = 87

# Scroll back so column A (not B) is the left-most visible column again,
# and leave A2:A6 selected (covering the new query cell).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A2:A6").Select()
